$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds a date serial number that was bumped by one day
# (from 2023-09-08 to 2023-09-09) for every data row (rows 2-397).
$ws.Range("C2:C397").Value = 45178
